# Season 14, matchdays prepares
# Inserts a new participant "Муратов Игорь" before the current row 13
# (pushing the existing rows 13-16 down to 14-17) and appends a new
# participant "Шевчук Антон" as a new row 18 at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Shift rows 13-16 down to 14-17 (bottom-up so we never overwrite
#    data we still need to read).
for ($r = 16; $r -ge 13; $r--) {
    $src = $ws.Range("A" + $r + ":W" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":W" + ($r + 1))
    $dst.Value2 = $src.Value2
}

# 2) The newly created row 17 needs the same formatting as the other
#    ranking cells in column A (bold, bordered, centered) - copy it
#    over from a neighbouring cell that already carries that style.
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# 3) Row 18 is a brand-new entrant with only a name, no results yet.
#    Write this one first so its shared-string gets allocated before
#    the row 13 entrant below (matches original authoring order).
$ws.Range("B18").Value = "Шевчук Антон"

# 4) Row 13 becomes the new entrant "Муратов Игорь" - only the ranking
#    cell keeps its style (value blank) and the name is filled in; all
#    of the old match data that used to live here has moved to row 14.
$ws.Range("C13:W13").ClearContents()
$ws.Range("A13").ClearContents()
$ws.Range("B13").Value = "Муратов Игорь"

# 5) Restore the view's active selection to match the latest edit.
$ws.Range("F13").Select() | Out-Null

# 6) Keep the worksheet's remembered sort range in sync with the table's
#    new extent.
$sort = $ws.Sort
$sort.SortFields.Clear() | Out-Null
$sort.SortFields.Add($ws.Range("B1:B17")) | Out-Null
$sort.SetRange($ws.Range("A2:W17"))
$sort.Header = 2
$sort.Apply()
